# Publish Latest checklists 2026-02-02
# Updates based on OWASP/wstg@e08e402
#
# Inserts a new "WSTG-CLNT-15 - Testing for Client-side Template Injection"
# row into the "Testing Checklist" sheet, just before the "API Testing"
# section (which pushes every row from the old row 133 onward down by one).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Testing Checklist")

# ---------------------------------------------------------------------
# 1) Insert a new blank row at 133 - shifts rows 133..138 down to 134..139
#    (data, formulas, styles and row heights all travel with the shift).
# ---------------------------------------------------------------------
$ws.Rows.Item(133).Insert()

# ---------------------------------------------------------------------
# 2) Populate the new row 133 with the WSTG-CLNT-15 checklist entry.
#    Styles mirror the neighbouring WSTG-CLNT-* rows (e.g. row 132).
# ---------------------------------------------------------------------
$ws.Rows.Item(133).RowHeight = 99

$ws.Cells.Item(133, 1).Style = $ws.Cells.Item(132, 1).Style
$ws.Cells.Item(133, 2).Style = $ws.Cells.Item(132, 2).Style
$ws.Cells.Item(133, 3).Style = $ws.Cells.Item(132, 3).Style
$ws.Cells.Item(133, 4).Style = $ws.Cells.Item(132, 4).Style
$ws.Cells.Item(133, 5).Style = $ws.Cells.Item(132, 5).Style
$ws.Cells.Item(133, 6).Style = $ws.Cells.Item(132, 6).Style

$ws.Range("B133").Value = "WSTG-CLNT-15"
$ws.Range("C133").Formula = '=HYPERLINK("https://owasp.org/www-project-web-security-testing-guide/latest/4-Web_Application_Security_Testing/11-Client-side_Testing/15-Testing_for_Client-Side_Template_Injection", "Testing for Client-side Template Injection")'
$ws.Range("D133").Value = "- Identify the client-side framework and its version used by the application.`n- Detect injection points where user input is reflected into the DOM and processed by the template engine.`n- Assess if the injection allows for arbitrary JavaScript execution (XSS) via the template syntax."
$ws.Range("E133").Value = "Not Started"

# ---------------------------------------------------------------------
# 3) Extend the conditional formatting range that used to stop at row 138
#    so that it now also covers the newly shifted-in row 139.
# ---------------------------------------------------------------------
$conditions = $ws.Range("B4:F139").FormatConditions
for ($i = 1; $i -le $conditions.Count; $i++) {
    $cond = $conditions.Item($i)
    if ($cond.AppliesTo.Address() -eq "$B$4:$F$138") {
        $cond.ModifyAppliesToRange($ws.Range("B4:F139"))
    }
}

# ---------------------------------------------------------------------
# 4) Data validation: the shared "Not Started,Pass,Issues,N/A" dropdown
#    rule needs to include the new E133 along with the shifted E136:E138
#    (the shift already renumbered E135/E136/E137 -> E136/E137/E138; we
#    only need to add the brand-new E133 cell to the rule).
# ---------------------------------------------------------------------
$ws.Range("E133").Validation.Delete()
$ws.Range("E133").Validation.Add(3, 1, 1, "Not Started,Pass,Issues,N/A")
